# Rename the repeated "Converting from Decimal to Binary" slide titles to
# "Converting from Binary to Decimal" across the whole deck.
$p = $ppt.ActivePresentation

$oldTitle = "Converting from Decimal to Binary"
$newTitle = "Converting from Binary to Decimal"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)

    if ($slide.Shapes.HasTitle) {
        $titleShape = $slide.Shapes.Title
        if ($titleShape.HasTextFrame) {
            if ($titleShape.TextFrame.TextRange.Text -eq $oldTitle) {
                $titleShape.TextFrame.TextRange.Text = $newTitle
            }
        }
    }
}
